# Add "2022-Q3" quarterly data: a new worksheet with the fund holdings for
# that quarter, plus a new leading row on the "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计" (i.e. as the new
#    2nd sheet), pushing every other quarter sheet one slot to the right.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# NOTE: look this sheet reference up only *after* the insert above - sheet
# handles obtained by name/index before an insert can end up pointing at the
# newly-inserted sheet once the tab order shifts.
$refSheet = $wb.Worksheets.Item("2022-Q2")

# Reuse the existing header / first-column styling (bold, centered, bordered)
# from the 2022-Q2 sheet instead of inventing new style entries.
$refSheet.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$refSheet.Range("A2").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q3Data = @(
    @(0, "009010", "华夏兴阳一年持有期混合",     "26.58", "88.59", "3.63", "0.9649", 4),
    @(1, "590008", "中邮战略新兴产业混合",       "7.25",  "90.68", "4.59", "0.3328", 5),
    @(2, "160425", "华安创业板两年定期开放混合", "1.80",  "93.72", "5.05", "0.0909", 8),
    @(3, "006992", "嘉合锦创优势精选混合",       "0.08",  "84.53", "2.89", "0.0023", 10)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = "'" + $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = "'" + $row[3]
    $q3.Range("E$r").Value = "'" + $row[4]
    $q3.Range("F$r").Value = "'" + $row[5]
    $q3.Range("G$r").Value = "'" + $row[6]
    $q3.Range("H$r").Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------------
# 2. Add the matching "2022-Q3" row at the top of the "总计" summary table.
# ---------------------------------------------------------------------------
$ws1 = $totalSheet

# Append a blank row 9 and give it the same formatting as the current last
# data row (8) so column A keeps its bold/centered/bordered look.
$ws1.Rows.Item(9).Insert()
$ws1.Range("A8").Copy()
$ws1.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Shift the existing quarters (rows 2-8) down one row (bottom-up so values
# aren't clobbered before they are read).
for ($row = 8; $row -ge 2; $row--) {
    $dst = $row + 1
    $ws1.Range("B$dst").Value = $ws1.Range("B$row").Value2
    $ws1.Range("C$dst").Value = $ws1.Range("C$row").Value2
    $ws1.Range("D$dst").Value = $ws1.Range("D$row").Value2
}

# Write the new 2022-Q3 summary figures into row 2.
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 1.39

# Re-sequence the leading index column (0..7) across rows 2..9.
for ($row = 2; $row -le 9; $row++) {
    $ws1.Range("A$row").Value = $row - 2
}

# ---------------------------------------------------------------------------
# 3. Restore the original active sheet/selection (总计, cell A1).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
